$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 66

# Copy formatting from the row above (row 65) down to the new row so that
# styles (bold/border on column A, date format on column E) are preserved.
$ws.Range("A65:V65").Copy() | Out-Null
$ws.Range("A66:V66").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = 65
$ws.Cells.Item($newRow, 2).Value = "azerbaijan"
$ws.Cells.Item($newRow, 3).Value = "premier-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45254.70833333334
$ws.Cells.Item($newRow, 6).Value = "Sabah Baku"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Qarabag"
$ws.Cells.Item($newRow, 9).Value = 2
$ws.Cells.Item($newRow, 10).Value = 3.64
$ws.Cells.Item($newRow, 11).Value = "23/11/2023 05:12"
$ws.Cells.Item($newRow, 12).Value = 4.54
$ws.Cells.Item($newRow, 13).Value = "24/11/2023 16:37"
$ws.Cells.Item($newRow, 14).Value = 3.31
$ws.Cells.Item($newRow, 15).Value = "23/11/2023 05:12"
$ws.Cells.Item($newRow, 16).Value = 3.58
$ws.Cells.Item($newRow, 17).Value = "24/11/2023 16:37"
$ws.Cells.Item($newRow, 18).Value = 1.89
$ws.Cells.Item($newRow, 19).Value = "23/11/2023 05:12"
$ws.Cells.Item($newRow, 20).Value = 1.75
$ws.Cells.Item($newRow, 21).Value = "24/11/2023 16:37"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/sabah-baku-qarabag-agdam/xp2vrm5c/"
